$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete entries at the top of the list ("Balite" and
# "Balungao"); this shifts every subsequent row up by two and naturally
# produces the new A1:H26 dimension.
$ws.Rows("2:3").Delete()

# The former "Buguion" row (now row 2) was renamed.
$ws.Range("B2").Value = "dicsa"

# A handful of entries flipped from Active back to not-Active as part of
# the dashboard / solveSettings refresh.
$ws.Range("A10").Value = $false
$ws.Range("A11").Value = $false
$ws.Range("A12").Value = $false
$ws.Range("A18").Value = $false
